$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (rows 1-4, columns A and B)
$ws.Range("A1").Value = -0.012450511837000217
$ws.Range("B1").Value = 0.012450511549646057

$ws.Range("A2").Value = -0.0069591349064320887
$ws.Range("B2").Value = 0.0069591346099887033

$ws.Range("A3").Value = -0.041315594993300155
$ws.Range("B3").Value = 0.041315594699510703

$ws.Range("A4").Value = 0.051833448318476472
$ws.Range("B4").Value = -0.051833448610635276

# Update column widths (col A widened, col B narrowed)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 12.833333333333334
